# REPORTGEN-141 : templates for TABLE_METRIC_ID_ROW
#
# Duplicate the "3-TableBlock-TableMetricIdCol" worksheet (placing the
# copy right after it) to create the new "3-TableBlock-TableMetricIdRow"
# template sheet, then update its three text cells (B1, B2, B12) to
# describe the new TABLE_METRIC_ID_ROW block instead of TABLE_METRIC_ID_COL.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("3-TableBlock-TableMetricIdCol")

# Copy the sheet, inserting the copy immediately after the source sheet.
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)

# The newly created copy becomes active and sits right after the source.
$newSheet = $wb.Worksheets.Item($srcSheet.Index + 1)
$newSheet.Name = "3-TableBlock-TableMetricIdRow"
$newSheet.Activate()

$newSheet.Range("B1").Value = "3.5. - TABLE_METRIC_ID_ROW"
$newSheet.Range("B2").Value = "* Block Name = TABLE_METRIC_ID_ROW"
$newSheet.Range("B12").Value = "RepGen:TABLE;TABLE_METRIC_ID_ROW;QID=60017|60014,SID=10151|67010,BID=66061,LEVEL=APPLICATION,SNAPSHOT=BOTH,VARIATION=BOTH,HEADER=SHORT"
